$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw")

# Replace the "ms" unit label with "us" across all the cells that used it.
# (the underlying shared-string table is rebuilt on save, so the now-unused
# "ms" entry is dropped automatically and "us" is appended/deduped)
$ws.Range("M1").Value2 = "us"
$ws.Range("M2").Value2 = "us"
$ws.Range("M3").Value2 = "us"
$ws.Range("M4").Value2 = "us"
$ws.Range("M5").Value2 = "us"
$ws.Range("J8").Value2 = "us"
$ws.Range("J9").Value2 = "us"
$ws.Range("J10").Value2 = "us"
$ws.Range("J11").Value2 = "us"
$ws.Range("J12").Value2 = "us"
$ws.Range("J13").Value2 = "us"

# New measurement row: period of 1080 (us) added to the L5 formula pattern.
$ws.Range("I6").Value2 = 1080
$ws.Range("L6").Formula = "=1000*1000/I6"

# Update the saved selection/active cell for the "raw" sheet.
$ws.Activate() | Out-Null
$ws.Range("J7").Select() | Out-Null
